# Update "want to go" counts (column F) for several rows across sheets,
# reflecting a refreshed scrape of show.bilibili.com data.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 12568
$ws1.Range("F10").Value = 980
$ws1.Range("F11").Value = 132
$ws1.Range("F12").Value = 339
$ws1.Range("F17").Value = 233
$ws1.Range("F19").Value = 17
$ws1.Range("F20").Value = 268
$ws1.Range("F21").Value = 297
$ws1.Range("F23").Value = 121
$ws1.Range("F25").Value = 5183
$ws1.Range("F27").Value = 1386
$ws1.Range("F29").Value = 1231

# --- Sheet "本地生活" (Local Life) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 9221
$ws3.Range("F4").Value = 1959

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 9221
$ws4.Range("F4").Value  = 1959
$ws4.Range("F5").Value  = 12568
$ws4.Range("F15").Value = 980
$ws4.Range("F16").Value = 132
$ws4.Range("F17").Value = 339
$ws4.Range("F22").Value = 233
$ws4.Range("F24").Value = 17
$ws4.Range("F25").Value = 268
$ws4.Range("F26").Value = 297
$ws4.Range("F33").Value = 5183
$ws4.Range("F35").Value = 1386
$ws4.Range("F40").Value = 1231
